$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.12%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'32.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.82%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.954"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.38%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07666"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.20%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.942"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-17.24%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.832"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.20%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.804"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.74%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1755"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.53%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07765"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.88%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08598"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-6.32%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03169"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.67%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.03%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005736"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.39%"
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'-0.22%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.153"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-4.17%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.3350"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.39%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1328"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.71%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.275"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.66%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.1995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'11.54%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04515"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.32%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001222"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.07%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.54%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001251"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.01690"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.50%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.45%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007465"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.25%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.67%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002332"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.59%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01046"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.31%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.90%"
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'0.19%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.8206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'10.43%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.19%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("E50").Style = "Normal"
